$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new header cells P1=14, Q1=15, using same style as existing header cells (s="1")
$ws.Cells.Item(1, 16).Value = 14   # P1
$ws.Cells.Item(1, 17).Value = 15   # Q1
$ws.Cells.Item(1, 15).Copy()
$ws.Range($ws.Cells.Item(1, 16), $ws.Cells.Item(1, 17)).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Rows 2-25: update I, K, M, O columns, and add P, Q columns
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I column: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K column: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M column: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O column: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P column: new
    $ws.Cells.Item($r, 17).Value = 2   # Q column: new
}
